$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Morning" was removed from the shared-strings table and a new
# "GIT UPDATE" entry was appended at the end; the only cell that
# referenced "Good Morning" (E8) now points at the new string.
$ws.Range("E8").Value = "GIT UPDATE"

# The edited file was last saved with E8 selected/active.
$ws.Range("E8").Select() | Out-Null
